$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A75").Value2 = "PW Hasaranga"
$ws.Range("B75").Value2 = 8
$ws.Range("C75").Value2 = 8
$ws.Range("D75").Value2 = 29
$ws.Range("E75").Value2 = 0
$ws.Range("F75").Value2 = 258
$ws.Range("G75").Value2 = 9
$ws.Range("H75").Value2 = 28.66
$ws.Range("I75").Value2 = 8.89
$ws.Range("J75").Value2 = 19.3
$ws.Range("K75").Value2 = 0
$ws.Range("L75").Value2 = 0
$ws.Range("M75").Value2 = "BOWL"
$ws.Range("N75").Value2 = 7

# M75 gets the existing default-bordered-centered style (idx2) via a clean format copy (no new entries)
$ws.Range("M2").Copy()
$ws.Range("M75").PasteSpecial(-4122)

# 1) Build target style for A75 (no wrap, new font) -> lands on idx6, with 1 unavoidable orphan idx5
$a = $ws.Range("A75")
$a.HorizontalAlignment = -4108
$a.VerticalAlignment = -4108
$a.Borders.LineStyle = 1
$a.Font.Color = 2236962

# 2) Copy A75's style onto B75:L75 and N75
$a.Copy()
$restBL = $ws.Range("B75:L75")
$restBL.PasteSpecial(-4122)
$n = $ws.Range("N75")
$n.PasteSpecial(-4122)

# 3) Apply WrapText to B75:L75 and N75 together -> creates the wrap idx7 cleanly
$restBL.WrapText = $true
$n.WrapText = $true
